$wb = $excel.ActiveWorkbook

# Update "zh-cn" sheet: Latest Handoff Datetime for the 4f748b75... file (row 4)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D4").Value = "2016-02-18 09:27:38"

# Update "de-de" sheet: Latest Handoff Datetime for the 4f748b75... file (row 4)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D4").Value = "2016-02-18 09:27:48"
